# Generate Report for Handback
# Adds a new handback record (730c43ea-7b8b-436c-b1ff-838e0662ebf5.md) as row 4
# to the Overview, zh-cn and de-de tables.

$wb = $excel.ActiveWorkbook

$fileName     = "730c43ea-7b8b-436c-b1ff-838e0662ebf5.md"
$pathAndName  = "e2e\730c43ea-7b8b-436c-b1ff-838e0662ebf5.md"
$commitSha    = "785c7eff2a9e5c08cbcdf2647011e25b545293f5"
$zhXlf        = "730c43ea-7b8b-436c-b1ff-838e0662ebf5.$commitSha.zh-cn.xlf"
$deXlf        = "730c43ea-7b8b-436c-b1ff-838e0662ebf5.$commitSha.de-de.xlf"

# ---------------------------------------------------------------------------
# Sheet "Overview" (table3 / Overview)
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item("Overview")
$loOverview.ListRows.Add() | Out-Null

$wsOverview.Range("A4").Value = $fileName
$wsOverview.Range("B4").Value = $pathAndName
$wsOverview.Range("C4").Value = ".md"
$wsOverview.Range("E4").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F4").Value = "Handed back: in sync with en-US"
$wsOverview.Range("G4").Value = "2016-08-25 20:44:08"

$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("B4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$commitSha/e2e/730c43ea-7b8b-436c-b1ff-838e0662ebf5.md",
    "",
    "",
    $pathAndName
) | Out-Null

# ---------------------------------------------------------------------------
# Sheet "zh-cn" (table1 / zh_cn)
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$loZh = $wsZh.ListObjects.Item("zh_cn")
$loZh.ListRows.Add() | Out-Null

$wsZh.Range("A4").Value = $fileName
$wsZh.Range("B4").Value = ".md"
$wsZh.Range("C4").Value = "Handed back: in sync with en-US"
$wsZh.Range("D4").Value = "e2e"
$wsZh.Range("E4").Value = "ht"
$wsZh.Range("F4").Value = "'True"
$wsZh.Range("G4").Value = $zhXlf
$wsZh.Range("H4").Value = "2016-08-25 20:43:57"
$wsZh.Range("I4").Value = $fileName
$wsZh.Range("J4").Value = $zhXlf
$wsZh.Range("K4").Value = "2016-08-25 20:44:28"
$wsZh.Range("L4").Value = ""
$wsZh.Range("M4").Value = "'True"
$wsZh.Range("N4").Value = ""
$wsZh.Range("O4").Value = "'False"
$wsZh.Range("P4").Value = ""

$wsZh.Hyperlinks.Add(
    $wsZh.Range("A4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$commitSha/e2e/730c43ea-7b8b-436c-b1ff-838e0662ebf5.md",
    "",
    "",
    $fileName
) | Out-Null

$wsZh.Hyperlinks.Add(
    $wsZh.Range("I4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/$commitSha/e2e/730c43ea-7b8b-436c-b1ff-838e0662ebf5.md",
    "",
    "",
    $fileName
) | Out-Null

# ---------------------------------------------------------------------------
# Sheet "de-de" (table2 / de_de)
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$loDe = $wsDe.ListObjects.Item("de_de")
$loDe.ListRows.Add() | Out-Null

$wsDe.Range("A4").Value = $fileName
$wsDe.Range("B4").Value = ".md"
$wsDe.Range("C4").Value = "Handed back: in sync with en-US"
$wsDe.Range("D4").Value = "e2e"
$wsDe.Range("E4").Value = "ht"
$wsDe.Range("F4").Value = "'True"
$wsDe.Range("G4").Value = $deXlf
$wsDe.Range("H4").Value = "2016-08-25 20:44:08"
$wsDe.Range("I4").Value = $fileName
$wsDe.Range("J4").Value = $deXlf
$wsDe.Range("K4").Value = "2016-08-25 20:44:35"
$wsDe.Range("L4").Value = ""
$wsDe.Range("M4").Value = "'True"
$wsDe.Range("N4").Value = ""
$wsDe.Range("O4").Value = "'False"
$wsDe.Range("P4").Value = ""

$wsDe.Hyperlinks.Add(
    $wsDe.Range("A4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$commitSha/e2e/730c43ea-7b8b-436c-b1ff-838e0662ebf5.md",
    "",
    "",
    $fileName
) | Out-Null

$wsDe.Hyperlinks.Add(
    $wsDe.Range("I4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/$commitSha/e2e/730c43ea-7b8b-436c-b1ff-838e0662ebf5.md",
    "",
    "",
    $fileName
) | Out-Null

Write-Host "Handback row added to Overview, zh-cn and de-de sheets."
